$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.955.90'
$ws.Range('E2').Value = '  +1.40%  '
$ws.Range('D3').Value = '3.478.64'
$ws.Range('E3').Value = '  +1.74%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '415.00'
$ws.Range('E5').Value = '  +1.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '130.77'
$ws.Range('E6').Value = '  +0.99%  '
$ws.Range('E7').Value = '  -1.89%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.729'
$ws.Range('E9').Value = '  -1.38%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.154'
$ws.Range('E10').Value = '  +6.86%  '
$ws.Range('E11').Value = '  -2.34%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '9.84'
$ws.Range('E12').Value = '  +4.46%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000224'
$ws.Range('E13').Value = '  -2.64%  '
$ws.Range('D14').Value = '4.027.59'
$ws.Range('E14').Value = '  +1.84%  '
$ws.Range('E15').Value = '  -0.35%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.51'
$ws.Range('E16').Value = '  -4.36%  '
$ws.Range('D17').Value = '3.479.77'
$ws.Range('E17').Value = '  +2.11%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.69'
$ws.Range('E18').Value = '  +0.95%  '
$ws.Range('E19').Value = '  +0.67%  '
$ws.Range('D20').Value = '62.840.69'
$ws.Range('E20').Value = '  +1.35%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '470.85'
$ws.Range('E21').Value = '  -1.94%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '90.81'
$ws.Range('E22').Value = '  -2.82%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.32'
$ws.Range('E23').Value = '  +3.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.20'
$ws.Range('E24').Value = '  -0.40%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.51'
$ws.Range('E25').Value = '  +12.36%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.34'
$ws.Range('E26').Value = '  +0.44%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '33.54'
$ws.Range('E27').Value = '  -0.46%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '4.80'
$ws.Range('E28').Value = '  +0.60%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.56'
$ws.Range('E29').Value = '  -1.27%  '
$ws.Range('E30').Value = '  +0.45%  '
$ws.Range('E31').Value = '  -1.07%  '
$ws.Range('E32').Value = '  -1.23%  '
$ws.Range('E33').Value = '  -1.55%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '41.00'
$ws.Range('E34').Value = '  -3.73%  '
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '58.08'
$ws.Range('E36').Value = '  +7.69%  '
$ws.Range('E37').Value = '  -3.54%  '
$ws.Range('E38').Value = '  +0.16%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.81'
$ws.Range('E39').Value = '  +8.04%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.05'
$ws.Range('E40').Value = '  +2.92%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '150.30'
$ws.Range('E41').Value = '  +4.14%  '
$ws.Range('E42').Value = '  -1.28%  '
$ws.Range('E43').Value = '  +0.97%  '
$ws.Range('B44').Value = 'NEARProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.45'
$ws.Range('E44').Value = '  +1.34%  '
$ws.Range('B45').Value = 'LidoDAOToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.33'
$ws.Range('E45').Value = '  -2.51%  '
$ws.Range('E46').Value = '  +2.01%  '
$ws.Range('D47').Value = '0.0₃0583'
$ws.Range('E47').Value = '  +31.45%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.40'
$ws.Range('E48').Value = '  +11.48%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '16.45'
$ws.Range('E49').Value = '  -1.61%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '22.22'
$ws.Range('E50').Value = '  -1.93%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.142'
$ws.Range('E51').Value = '  -5.07%  '
